$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Suraj Kumar Jha"
$ws.Range("B3").Value = [datetime]"2000-09-19"

# The new row's text wraps onto two lines (same as row 2 for "Vikram Kumar
# Jha"), so Excel auto-grows the row height to fit the wrapped text.
$ws.Rows.Item(3).RowHeight = 29

$ws.Range("A4").Select()
